$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.761.16'
$ws.Range("E2").Value = '  +3.91%  '
$ws.Range("D3").Value = '2.270.15'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.01'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.22'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +5.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.530'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +2.86%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.40'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +4.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.26'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.60'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").Value = '2.619.22'
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.23'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +2.64%  '
$ws.Range("D17").Value = '2.288.79'
$ws.Range("E17").Value = '  +4.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.764'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +3.70%  '
$ws.Range("D19").Value = '41.687.40'
$ws.Range("E19").Value = '  +3.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.58'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +10.60%  '
$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.87'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.43'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("E25").Value = '  +3.42%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +5.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.08'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.52'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.07'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.71'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +2.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.18'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +5.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.25'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +5.40%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0744'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +3.83%  '
$ws.Range("E36").Value = '  +1.37%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.87'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +8.66%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("E41").Value = '  +3.93%  '
$ws.Range("E42").Value = '  +3.17%  '
$ws.Range("D43").Value = '2.047.84'
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.42'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0279'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.07'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +11.33%  '
$ws.Range("E48").Value = '  +2.17%  '
$ws.Range("E49").Value = '  +4.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.88'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +6.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.16'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +2.35%  '
